$d = $word.ActiveDocument

# Title: "Testing" + " " + "custom" + " " + "properties" -> "Testing custom properties"
$d.Content.Find.Execute("Testing custom properties", $true, $false, $false, $false, $false, $true, 1, $false, "Testing custom properties", 2)

# Subtitle: "This" + " " + "is" + " " + "a" + " " + "subtitle" -> "This is a subtitle"
$d.Content.Find.Execute("This is a subtitle", $true, $false, $false, $false, $false, $true, 1, $false, "This is a subtitle", 2)

# Author: "A." + " " + "M." -> "A. M."
$d.Content.Find.Execute("A. M.", $true, $false, $false, $false, $false, $true, 1, $false, "A. M.", 2)

# Abstract body: merge "Quite a long description" into one run, and
# "spanning several lines" into another, keeping the middle space run intact.
$d.Content.Find.Execute("Quite a long description", $true, $false, $false, $false, $false, $true, 1, $false, "Quite a long description", 2)
$d.Content.Find.Execute("spanning several lines", $true, $false, $false, $false, $false, $true, 1, $false, "spanning several lines", 2)
